{"js": "// Fix typo in User Guide: \"To login, click the Login link...\" should read\n// \"To log in, click the Login link...\" \u2014 i.e. split the single run\n// containing the word \"login\" into three runs: \"log\", \" \", \"in\".\n\nconst doc = context.document;\nconst body = doc.body;\n\n// Locate the specific sentence \"To login,\" (there are other, unrelated\n// occurrences of the word \"login\" elsewhere in the document \u2014 e.g.\n// \"...login there.\" and \"...your login information.\" \u2014 which must NOT\n// be touched).\nconst paraResults = body.search(\"To login,\", { matchCase: true });\nparaResults.load(\"items\");\nawait context.sync();\n\nif (paraResults.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly 1 match for 'To login,' but found \" + paraResults.items.length\n  );\n}\nconst hit = paraResults.items[0];\n\n// Within that hit, find \"log\" (the first 3 letters of \"login\") so we can\n// get the insertion point right between \"log\" and \"in\".\nconst logResults = hit.search(\"log\", { matchCase: true });\nlogResults.load(\"items\");\nawait context.sync();\n\nif (logResults.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly 1 match for 'log' inside 'To login,' but found \" + logResults.items.length\n  );\n}\nconst logRange = logResults.items[0];\nconst insertionPoint = logRange.getRange(\"After\");\nawait context.sync();\n\n// Plain insertText() here would silently re-merge the new text into the\n// neighboring run (ending up with a single run \"log in\" instead of three\n// runs \"log\" / \" \" / \"in\"). Turning change tracking on first causes the\n// inserted space to be recorded as its own run; accepting the revision\n// afterwards drops the <w:ins> wrapper but leaves the run split intact,\n// matching the target run structure.\ndoc.changeTrackingMode = \"TrackAll\";\nawait context.sync();\n\ninsertionPoint.insertText(\" \", \"Before\");\nawait context.sync();\n\ndoc.changeTrackingMode = \"Off\";\nawait context.sync();\n\nbody.getTrackedChanges().acceptAll();\nawait context.sync();\n", "ps1": "# Fix typo in User Guide: \"To login, click the Login link...\" should read\n# \"To log in, click the Login link...\" \u2014 i.e. split the single run\n# containing the word \"login\" into three runs: \"log\", \" \", \"in\".\n\n$d = $word.ActiveDocument\n\n# Locate the specific paragraph that starts with \"To login,\" (there are\n# other, unrelated occurrences of the word \"login\" elsewhere in the\n# document \u2014 e.g. \"...login there.\" and \"...your login information.\" \u2014\n# which must NOT be touched).\n$targetPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs($i)\n    if ($p.Range.Text -like \"To login,*\") {\n        $targetPara = $p\n        break\n    }\n}\n\nif ($targetPara -eq $null) {\n    throw \"Could not find the target paragraph ('To login, ...').\"\n}\n\n# Search for the standalone word \"login\" within just that paragraph.\n$rng = $targetPara.Range\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"login\"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $true\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\n$found = $rng.Find.Execute()\n\nif (-not $found) {\n    throw \"Could not find the word 'login' in the target paragraph.\"\n}\n\n# $rng now spans exactly the 5 characters \"login\". Insert a space after\n# the 3rd character (\"log\" | \"in\") while Track Changes is on: the tracked\n# insertion is recorded as its own run, which keeps \"log\" / \" \" / \"in\" as\n# three separate sibling runs (matching the target XML) instead of the\n# plain-insert behavior of silently re-merging into one run. Accepting\n# the revision afterwards then drops the <w:ins> wrapper but leaves the\n# run split intact.\n$wasTracking = $d.TrackRevisions\n$d.TrackRevisions = $true\n\n$splitPoint = $rng.Start + 3\n$splitRange = $d.Range($splitPoint, $splitPoint)\n$splitRange.InsertAfter(\" \")\n\n$d.TrackRevisions = $wasTracking\n$d.AcceptAllRevisions()\n"}
